$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 14 (M564), pushing it down to row 17.
# This keeps the existing row's data intact while making room for the three new
# cohort-9 animals (F523, F524, F525) to occupy rows 14-16.
$ws.Rows.Item(14).Resize(2).EntireRow.Insert()

# Row 14: F523
$ws.Cells.Item(14, 1).Value = "C05"
$ws.Cells.Item(14, 2).Value = "F523"
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = "K2C05HSSHA10-20180810.txt"
$ws.Cells.Item(14, 5).Value = "933000120138754"
$ws.Cells.Item(14, 6).Value = "F"
$ws.Cells.Item(14, 7).Value = 10
$ws.Cells.Item(14, 8).Value = "fail"
$ws.Cells.Item(14, 11).Value = 10

# Row 15: F524
$ws.Cells.Item(15, 1).Value = "C05"
$ws.Cells.Item(15, 2).Value = "F524"
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = "K2C05HSSHA10-20180810.txt"
$ws.Cells.Item(15, 5).Value = "933000120138743"
$ws.Cells.Item(15, 6).Value = "F"
$ws.Cells.Item(15, 7).Value = 2
$ws.Cells.Item(15, 8).Value = "fail"
$ws.Cells.Item(15, 11).Value = 2

# Row 16: F525
$ws.Cells.Item(16, 1).Value = "C05"
$ws.Cells.Item(16, 2).Value = "F525"
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = "K2C05HSSHA10-20180810.txt"
$ws.Cells.Item(16, 5).Value = "933000120138732"
$ws.Cells.Item(16, 6).Value = "F"
$ws.Cells.Item(16, 7).Value = 41
$ws.Cells.Item(16, 8).Value = "fail"
$ws.Cells.Item(16, 11).Value = 41

# Row 17 (the original row 14, now shifted down) gets its SHA10 count filled in.
$ws.Cells.Item(17, 11).Value = 10
